# Generate Report for Handoff
# The 42ef5cfb-...-md source file has finished translation prep and is now
# "Ready for handoff" — update its Status (and the Overview roll-up) plus
# the corresponding handoff timestamps for both locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 42ef5cfb-...md file
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 20:22:10"

# zh-cn sheet: row 3 (Status + Latest Handoff Datetime)
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-24 20:22:05"

# de-de sheet: row 3 (Status + Latest Handoff Datetime)
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-24 20:22:10"
